# Insert a new "ConstructionYear_1918" column between the existing
# ConstructionYear_1900 (B) and ConstructionYear_1945 (old C, now D) columns,
# then (re)write the full data table to match the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:M -> D:N, creating a blank column C for the new year bucket.
$ws.Columns("C").Insert()

# ---- Header row ----
$ws.Range("A1").Value = "Water"
$ws.Range("B1").Value = "ConstructionYear_1900"
$ws.Range("C1").Value = "ConstructionYear_1918"
$ws.Range("D1").Value = "ConstructionYear_1945"
$ws.Range("E1").Value = "ConstructionYear_1961"
$ws.Range("F1").Value = "ConstructionYear_1970"
$ws.Range("G1").Value = "ConstructionYear_1981"
$ws.Range("H1").Value = "ConstructionYear_1991"
$ws.Range("I1").Value = "ConstructionYear_2001"
$ws.Range("J1").Value = "ConstructionYear_2011"
$ws.Range("K1").Value = "ConstructionYear_2016"
$ws.Range("L1").Value = "ConstructionYear_2020"
$ws.Range("M1").Value = "ConstructionYear_2030"
$ws.Range("N1").Value = "ConstructionYear_2035"

# ---- Row 2: Low ----
$ws.Range("A2").Value = "Low"
$ws.Range("B2").Value = 756
$ws.Range("C2").Value = 2107
$ws.Range("D2").Value = 31728
$ws.Range("E2").Value = 24695
$ws.Range("F2").Value = 20596
$ws.Range("G2").Value = 19762
$ws.Range("H2").Value = 41657
$ws.Range("I2").Value = 36786
$ws.Range("J2").Value = 37786
$ws.Range("K2").Value = 15789
$ws.Range("L2").Value = 17284
$ws.Range("M2").Value = 27518
$ws.Range("N2").Value = 27313.31

# ---- Row 3: mid-Low ----
$ws.Range("A3").Value = "mid-Low"
$ws.Range("B3").Value = 3783
$ws.Range("D3").Value = 23941.59
$ws.Range("E3").Value = 33933.6
$ws.Range("F3").Value = 12107.27
$ws.Range("G3").Value = 21273.51
$ws.Range("H3").Value = 27826.15
$ws.Range("I3").Value = 25649.13
$ws.Range("J3").Value = 36460.21
$ws.Range("K3").Value = 15256.25
$ws.Range("L3").Value = 11792.58
$ws.Range("M3").Value = 25899.59
$ws.Range("N3").Value = 14275.55

# ---- Row 4: Middle ----
$ws.Range("A4").Value = "Middle"
$ws.Range("B4").Value = 5432.62
$ws.Range("C4").Value = 38.32
$ws.Range("D4").Value = 15585.17
$ws.Range("E4").Value = 21938.42
$ws.Range("F4").Value = 11685.87
$ws.Range("G4").Value = 25276.22
$ws.Range("H4").Value = 16113.1
$ws.Range("I4").Value = 30999.13
$ws.Range("J4").Value = 38580.44
$ws.Range("K4").Value = 26178.14
$ws.Range("L4").Value = 16685.87
$ws.Range("M4").Value = 32066.15
$ws.Range("N4").Value = 15668.16

# ---- Row 5: mid-High ----
$ws.Range("A5").Value = "mid-High"
$ws.Range("B5").Value = 910.64
$ws.Range("D5").Value = 17586.04
$ws.Range("E5").Value = 23999.04
$ws.Range("F5").Value = 10350.21
$ws.Range("G5").Value = 20567.05
$ws.Range("H5").Value = 33363.96
$ws.Range("I5").Value = 31628.68
$ws.Range("J5").Value = 19203.35
$ws.Range("K5").Value = 14848.76
$ws.Range("L5").Value = 12146.85
$ws.Range("M5").Value = 18881.18
$ws.Range("N5").Value = 10398.89

# ---- Row 6: High ----
$ws.Range("A6").Value = "High"
$ws.Range("B6").Value = 1989.63
$ws.Range("D6").Value = 14134.4
$ws.Range("E6").Value = 13426.42
$ws.Range("F6").Value = 20346.68
$ws.Range("G6").Value = 22532.3
$ws.Range("H6").Value = 24775.92
$ws.Range("I6").Value = 35835.83
$ws.Range("J6").Value = 31014.19
$ws.Range("K6").Value = 13740.42
$ws.Range("L6").Value = 6450.16
$ws.Range("M6").Value = 24354.21
$ws.Range("N6").Value = 7430.19
